# "Add files via upload" — refresh Sheet1's small demo grid and UI state
# to match the newly uploaded version of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Shift the block-pattern values in columns D:I for rows 2 and 3 down by
# one "block" (each pair of columns holds the same value), e.g.
# D2:E2 1,1 -> 0,0 ; F2:G2 2,2 -> 1,1 ; H2:I2 3,3 -> 2,2 (same for row 3).
$ws.Range("D2:E2").Value = 0
$ws.Range("F2:G2").Value = 1
$ws.Range("H2:I2").Value = 2

$ws.Range("D3:E3").Value = 0
$ws.Range("F3:G3").Value = 1
$ws.Range("H3:I3").Value = 2

# The author's selection moved from E14 to K9 before saving.
$ws.Range("K9").Select()

# The Excel window was resized wider before saving (cosmetic UI state).
$excel.ActiveWindow.Width = 15684
